# Adds a new "2022" data column (column T) to the table, mirroring the
# formatting of the existing column S, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the 2022 column (row -> value)
$newValues = @{
    4  = 2022
    5  = 3.7
    6  = 1.6
    7  = 1.7
    8  = 17.9
    9  = 7.5
    10 = 1.1
    11 = 4.4
    12 = 3
    13 = 4.1
    14 = 0.8
}

foreach ($row in 4..14) {
    # Copy formatting (number format, font, borders, alignment, etc.) from
    # the last existing year column (S) into the new column (T).
    $ws.Range("S$row").Copy($ws.Range("T$row"))
    $ws.Range("T$row").Value = $newValues[$row]
}

# Move / update the active selection as recorded in the workbook.
$ws.Range("U4").Select()
